$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new weekly record is inserted as row 198, pushing the former
# rows 198-289 down to 199-290 (a pure shift, no other data changes).
$ws.Rows.Item(198).Insert()

$ws.Cells.Item(198, 1).Value = 8
$ws.Cells.Item(198, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(198, 3).Value = "Coquimbo"
$ws.Cells.Item(198, 4).Value = 44609
$ws.Cells.Item(198, 5).Value = 4
$ws.Cells.Item(198, 6).Value = 100114013
$ws.Cells.Item(198, 7).Value = "Zanahoria"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 520
$ws.Cells.Item(198, 11).Value = 5800
$ws.Cells.Item(198, 12).Value = 6000
$ws.Cells.Item(198, 13).Value = 5900
$ws.Cells.Item(198, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(198, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(198, 16).Value = 295
$ws.Cells.Item(198, 17).Value = 20
$ws.Cells.Item(198, 18).Value = "Hortaliza"
